$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the disclaimer text date from 2021-05-28 to 2021-06-09
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-06-09 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-7
$ws.Range("D2").Value = 0.2499796881609357
$ws.Range("E2").Value = -0.0004887585532746819

$ws.Range("D3").Value = 0.5400052123275155
$ws.Range("E3").Value = -0.004445606694560733

$ws.Range("D4").Value = 0.0500129765665825
$ws.Range("E4").Value = -0.002895193977996535

$ws.Range("D5").Value = 0.1000047012343144
$ws.Range("E5").Value = -0.007232788641843069

$ws.Range("D6").Value = 0.05999742171065175
$ws.Range("E6").Value = -0.006524575902566276

$ws.Range("D7").Value = 0.9999999999999998
$ws.Range("E7").Value = -0.003782398365459794
